# Marksheet update: correct/total marks for the "Marking" and "Total" rows.
#   B11 (Marking, correct count)  : 3       -> 5
#   B12 (Total,   correct count)  : 45      -> 75
#   E12 (Total,   "correct/total"): 43/84   -> 75/140

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("quiz")
if (-not $ws) {
    $ws = $wb.ActiveSheet
}

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 75
$ws.Range("E12").Value = "75/140"
